# [Draft physical database] Add database diagram
# Adds a "Daily Meeting" attendance table to the "Impediment Backlog" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Impediment Backlog")

# --- New column D width/formatting setup ---
$ws.Columns.Item(4).ColumnWidth = 49.6

# Header row for the daily-meeting attendance block
$ws.Range("A15").Value = "Daily Meeting"
$ws.Range("C15").Value = "Absent Without Plans"
$ws.Range("D15").Value = "Absent With Plans"

# Date column (A16:A31) uses a dd/mm number format, but the values are kept
# as plain text (matches author's literal "13/05"-style strings).
$dateCells = @(
    @{ Row = 16; Date = "13/05" },
    @{ Row = 17; Date = "14/5" },
    @{ Row = 18; Date = "15/5" },
    @{ Row = 19; Date = "16/5" },
    @{ Row = 20; Date = "17/5" },
    @{ Row = 21; Date = "18/5" },
    @{ Row = 22; Date = "19/05" },
    @{ Row = 23; Date = "20/05" },
    @{ Row = 24; Date = "21/05" },
    @{ Row = 25; Date = "22/5" },
    @{ Row = 26; Date = "23/5" },
    @{ Row = 27; Date = "24/5" },
    @{ Row = 28; Date = "25/5" },
    @{ Row = 29; Date = "26/5" },
    @{ Row = 30; Date = "27/05" },
    @{ Row = 31; Date = "28/05" }
)

foreach ($entry in $dateCells) {
    $cell = $ws.Cells.Item($entry.Row, 1)
    $cell.NumberFormat = "dd/mm"
    $cell.Value = $entry.Date
}

# Notes / attendance remarks in columns C and D for specific days
$ws.Range("C16").Value = "Tín"
$ws.Range("C17").Value = "Tin"
$ws.Range("C20").Value = "Dang"
$ws.Range("D22").Value = "Tin"
$ws.Range("C23").Value = "Tin"
$ws.Range("C24").Value = "Tin"
$ws.Range("D28").Value = "Tin"
$ws.Range("D29").Value = "Dang"

# Restore view state: scroll so row 14 is near the top, with D28 selected
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$ws.Range("D28").Select() | Out-Null
